$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -1.85663703414174
$ws.Range("E2").Value = 0.0464444315132516
$ws.Range("F2").Value = 2.55591184998261
$ws.Range("G2").Value = 0.0910253476343038
$ws.Range("H2").Value = 0.807658553218601
